{"js": "// Change \"Two native species of flax...\" to \"The two native species of flax...\"\n// i.e. prepend \"The \" and lower-case the original \"Two\" -> \"two\".\nconst searchResults = context.document.body.search(\"Two native species of flax\", { matchCase: true });\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  searchResults.items[0].insertText(\n    \"The two native species of flax\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "# Change \"Two native species of flax...\" to \"The two native species of flax...\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n# wdFindContinue = 1, wdReplaceOne = 1\n$find.Execute(\n    \"Two native species of flax\",\n    $true,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"The two native species of flax\",\n    1\n) | Out-Null\n"}
